# Update cryptos list data (Price and Volume(1h) columns) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure D and E columns are formatted as Text so numeric-looking strings
# like "1.00" or "55.861.26" are preserved exactly as text, not coerced to numbers/dates.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "55.861.26"
$ws.Cells.Item(2, 5).Value = "  -2.19%  "
$ws.Cells.Item(3, 4).Value = "2.973.06"
$ws.Cells.Item(3, 5).Value = "  -0.53%  "
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
$ws.Cells.Item(5, 4).Value = "502.23"
$ws.Cells.Item(5, 5).Value = "  -0.08%  "
$ws.Cells.Item(6, 4).Value = "136.79"
$ws.Cells.Item(6, 5).Value = "  -1.31%  "
$ws.Cells.Item(7, 5).Value = "  +0.12%  "
$ws.Cells.Item(8, 4).Value = "0.425"
$ws.Cells.Item(8, 5).Value = "  -1.31%  "
$ws.Cells.Item(9, 5).Value = "  -2.55%  "
$ws.Cells.Item(10, 5).Value = "  -1.88%  "
$ws.Cells.Item(11, 4).Value = "0.362"
$ws.Cells.Item(11, 5).Value = "  +0.58%  "
$ws.Cells.Item(12, 4).Value = "3.481.90"
$ws.Cells.Item(12, 5).Value = "  -0.28%  "
$ws.Cells.Item(13, 5).Value = "  -1.66%  "
$ws.Cells.Item(14, 4).Value = "25.87"
$ws.Cells.Item(14, 5).Value = "  -0.88%  "
$ws.Cells.Item(15, 5).Value = "  -0.65%  "
$ws.Cells.Item(16, 4).Value = "55.876.98"
$ws.Cells.Item(16, 5).Value = "  -2.18%  "
$ws.Cells.Item(17, 4).Value = "2.971.27"
$ws.Cells.Item(17, 5).Value = "  -0.38%  "
$ws.Cells.Item(18, 4).Value = "5.98"
$ws.Cells.Item(18, 5).Value = "  -1.53%  "
$ws.Cells.Item(19, 4).Value = "12.83"
$ws.Cells.Item(19, 5).Value = "  +1.25%  "
$ws.Cells.Item(20, 4).Value = "7.95"
$ws.Cells.Item(20, 5).Value = "  +0.66%  "
$ws.Cells.Item(21, 4).Value = "326.87"
$ws.Cells.Item(21, 5).Value = "  +1.55%  "
$ws.Cells.Item(22, 5).Value = "  +0.22%  "
$ws.Cells.Item(23, 4).Value = "0.491"
$ws.Cells.Item(23, 5).Value = "  +0.07%  "
$ws.Cells.Item(24, 4).Value = "64.29"
$ws.Cells.Item(24, 5).Value = "  +0.64%  "
$ws.Cells.Item(25, 4).Value = "3.091.04"
$ws.Cells.Item(25, 5).Value = "  -0.43%  "
$ws.Cells.Item(26, 4).Value = "0.999"
$ws.Cells.Item(26, 5).Value = "  +0.04%  "
$ws.Cells.Item(27, 5).Value = "  -1.90%  "
$ws.Cells.Item(28, 4).Value = "0.0₃0889"
$ws.Cells.Item(28, 5).Value = "  -1.28%  "
$ws.Cells.Item(29, 4).Value = "6.35"
$ws.Cells.Item(29, 5).Value = "  -3.55%  "
$ws.Cells.Item(30, 4).Value = "6.93"
$ws.Cells.Item(30, 5).Value = "  -2.23%  "
$ws.Cells.Item(31, 4).Value = "1.77"
$ws.Cells.Item(31, 5).Value = "  -0.54%  "
$ws.Cells.Item(32, 4).Value = "20.11"
$ws.Cells.Item(32, 5).Value = "  -0.71%  "
$ws.Cells.Item(33, 5).Value = "  -2.29%  "
$ws.Cells.Item(34, 4).Value = "153.24"
$ws.Cells.Item(34, 5).Value = "  -1.44%  "
$ws.Cells.Item(35, 4).Value = "4.47"
$ws.Cells.Item(35, 5).Value = "  -2.71%  "
$ws.Cells.Item(36, 4).Value = "5.68"
$ws.Cells.Item(36, 5).Value = "  -2.13%  "
$ws.Cells.Item(37, 4).Value = "25.42"
$ws.Cells.Item(37, 5).Value = "  +5.12%  "
$ws.Cells.Item(38, 4).Value = "1.23"
$ws.Cells.Item(38, 5).Value = "  -2.16%  "
$ws.Cells.Item(39, 4).Value = "0.0654"
$ws.Cells.Item(39, 5).Value = "  -2.12%  "
$ws.Cells.Item(40, 4).Value = "3.006.25"
$ws.Cells.Item(40, 5).Value = "  -0.42%  "
$ws.Cells.Item(41, 4).Value = "36.73"
$ws.Cells.Item(41, 5).Value = "  -3.05%  "
$ws.Cells.Item(42, 4).Value = "1.00"
$ws.Cells.Item(42, 5).Value = "  +0.09%  "
$ws.Cells.Item(43, 4).Value = "3.76"
$ws.Cells.Item(43, 5).Value = "  -0.11%  "
$ws.Cells.Item(44, 4).Value = "0.647"
$ws.Cells.Item(44, 5).Value = "  +0.80%  "
$ws.Cells.Item(45, 4).Value = "2.152.55"
$ws.Cells.Item(45, 5).Value = "  -2.31%  "
$ws.Cells.Item(46, 5).Value = "  -3.80%  "
$ws.Cells.Item(47, 4).Value = "5.80"
$ws.Cells.Item(47, 5).Value = "  -3.12%  "
$ws.Cells.Item(48, 4).Value = "0.917"
$ws.Cells.Item(48, 5).Value = "  -2.77%  "
$ws.Cells.Item(49, 4).Value = "0.0234"
$ws.Cells.Item(49, 5).Value = "  -0.99%  "
$ws.Cells.Item(50, 4).Value = "19.46"
$ws.Cells.Item(50, 5).Value = "  +0.43%  "
$ws.Cells.Item(51, 4).Value = "0.0846"
$ws.Cells.Item(51, 5).Value = "  -3.79%  "
